$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2

# 1) "...is used to reduce the result a Curried expression..."
#    -> "...is used to reduce the result of a Curried expression..."
$d.Content.Find.Execute(
    "reduce the result a Curried",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "reduce the result of a Curried",
    2)

# 2) "...with the name SK.  They appear in bold face..."
#    -> "...with the name SK.  These symbols appear in bold face..."
$d.Content.Find.Execute(
    "with the name SK.  They appear",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "with the name SK.  These symbols appear",
    2)

# 3) "...String and Numeric literals and operators will also be explained in a subsequent edition."
#    -> "...String and Numeric literals and their operators will be documented in a subsequent edition."
$d.Content.Find.Execute(
    "String and Numeric literals and operators will also be explained in a subsequent edition.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "String and Numeric literals and their operators will be documented in a subsequent edition.",
    2)
